$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1688.8695
$ws.Range("I15").Value = 1688.8695
$ws.Range("K15").Value = 5066.6085
$ws.Range("M15").Value = -4897.6085

$ws.Range("H17").Value = 1184650.4
$ws.Range("J17").Value = 1381921.2
$ws.Range("L17").Value = 4145763.6
$ws.Range("N17").Value = -4146099.6

$ws.Range("H19").Value = 1947.6364
$ws.Range("J19").Value = 393.625
$ws.Range("L19").Value = 393.625
$ws.Range("N19").Value = -743.625

$ws.Range("H21").Value = 999
$ws.Range("I21").Value = 999
$ws.Range("K21").Value = 999
$ws.Range("M21").Value = -531

$ws.Range("H23").Value = 999
$ws.Range("I23").Value = 999
$ws.Range("K23").Value = 999
$ws.Range("M23").Value = -765

$ws.Range("H33").Value = 785027.0600000001
$ws.Range("I33").Value = 1015110.6
$ws.Range("J33").Value = 2743
$ws.Range("K33").Value = 1015110.6
$ws.Range("L33").Value = 2743
$ws.Range("M33").Value = -1014881.6
$ws.Range("N33").Value = -3201

$ws.Range("H61").Value = 838.5
$ws.Range("I61").Value = 802.8
$ws.Range("J61").Value = 1017
$ws.Range("K61").Value = 2408.4
$ws.Range("L61").Value = 3051
$ws.Range("M61").Value = -2236.4
$ws.Range("N61").Value = -3395

$ws.Range("H74").Value = 4077.077
$ws.Range("I74").Value = 1875.5
$ws.Range("K74").Value = 1875.5
$ws.Range("M74").Value = -939.5

$ws.Range("H77").Value = 4077.077
$ws.Range("I77").Value = 1875.5
$ws.Range("K77").Value = 9377.5
$ws.Range("M77").Value = -4697.5

$ws.Range("H80").Value = 935.3333
$ws.Range("I80").Value = 780.3333
$ws.Range("J80").Value = 1038.6666
$ws.Range("K80").Value = 2340.9999
$ws.Range("L80").Value = 3115.9998
$ws.Range("M80").Value = -1342.9999
$ws.Range("N80").Value = -5111.9998

$ws.Range("H83").Value = 935.3333
$ws.Range("I83").Value = 780.3333
$ws.Range("J83").Value = 1038.6666
$ws.Range("K83").Value = 7022.9997
$ws.Range("L83").Value = 9347.999400000001
$ws.Range("M83").Value = -2030.9997
$ws.Range("N83").Value = -19331.9994

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H98").Value = 3776.2307
$ws.Range("I98").Value = 3281.913
$ws.Range("K98").Value = 3281.913
$ws.Range("M98").Value = -1783.913

$ws.Range("H106").Value = 12340.091
$ws.Range("I106").Value = 2962.5
$ws.Range("K106").Value = 2962.5
$ws.Range("M106").Value = -2331.5

$ws.Range("H113").Value = 4725.1113
$ws.Range("J113").Value = 5331.8335
$ws.Range("L113").Value = 5331.8335
$ws.Range("N113").Value = -11839.8335

$ws.Range("H116").Value = 7691.077
$ws.Range("I116").Value = 4997.5
$ws.Range("K116").Value = 4997.5
$ws.Range("M116").Value = -1555.5

$ws.Range("H122").Value = 3776.2307
$ws.Range("I122").Value = 3281.913
$ws.Range("K122").Value = 9845.739
$ws.Range("M122").Value = -7395.739

$ws.Range("H135").Value = 1078.4286
$ws.Range("I135").Value = 624.75
$ws.Range("K135").Value = 5622.75
$ws.Range("M135").Value = -3087.75

$ws.Range("H136").Value = 175780
$ws.Range("J136").Value = 175780
$ws.Range("L136").Value = 175780
$ws.Range("N136").Value = -185980

$ws.Range("H137").Value = 68552.13
$ws.Range("I137").Value = 2019.4
$ws.Range("J137").Value = 201617.6
$ws.Range("K137").Value = 6058.200000000001
$ws.Range("L137").Value = 604852.8
$ws.Range("M137").Value = -3508.200000000001
$ws.Range("N137").Value = -609952.8

$ws.Range("H138").Value = 3636.647
$ws.Range("J138").Value = 5139
$ws.Range("L138").Value = 15417
$ws.Range("N138").Value = -25697

$ws.Range("H141").Value = 47935.95
$ws.Range("I141").Value = 55742.06
$ws.Range("K141").Value = 167226.18
$ws.Range("M141").Value = -162046.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 252886.33
$ws.Range("I32").Value = 288354.34
$ws.Range("K32").Value = 288354.34
$ws.Range("M32").Value = -288067.34

$ws.Range("H48").Value = 125066
$ws.Range("J48").Value = 125066
$ws.Range("L48").Value = 125066
$ws.Range("N48").Value = -125834

$ws.Range("H61").Value = 3456.5833
$ws.Range("I61").Value = 3456.5833
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3456.5833
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3244.5833
$ws.Range("N61").ClearContents()

$ws.Range("H63").Value = 1150
$ws.Range("I63").Value = 1150
$ws.Range("K63").Value = 1150
$ws.Range("M63").Value = -464

$ws.Range("H66").Value = 1150
$ws.Range("I66").Value = 1150
$ws.Range("K66").Value = 5750
$ws.Range("M66").Value = -2318

$ws.Range("H74").Value = 1498.4642
$ws.Range("I74").Value = 957.3158
$ws.Range("J74").Value = 2640.889
$ws.Range("K74").Value = 957.3158
$ws.Range("L74").Value = 2640.889
$ws.Range("M74").Value = -83.31579999999997
$ws.Range("N74").Value = -4388.889

$ws.Range("H77").Value = 1498.4642
$ws.Range("I77").Value = 957.3158
$ws.Range("J77").Value = 2640.889
$ws.Range("K77").Value = 4786.579
$ws.Range("L77").Value = 13204.445
$ws.Range("M77").Value = -418.5789999999997
$ws.Range("N77").Value = -21940.445

$ws.Range("H95").Value = 163707.75
$ws.Range("J95").Value = 163707.75
$ws.Range("L95").Value = 163707.75
$ws.Range("N95").Value = -169199.75

$ws.Range("H96").Value = 79999
$ws.Range("J96").Value = 79999
$ws.Range("L96").Value = 79999
$ws.Range("N96").Value = -85491

$ws.Range("H122").Value = 10271.083
$ws.Range("I122").Value = 11776.3
$ws.Range("J122").Value = 2745
$ws.Range("K122").Value = 35328.89999999999
$ws.Range("L122").Value = 8235
$ws.Range("M122").Value = -32878.89999999999
$ws.Range("N122").Value = -13135

$ws.Range("H132").Value = 1082.3636
$ws.Range("I132").Value = 1082.3636
$ws.Range("K132").Value = 3247.0908
$ws.Range("M132").Value = -717.0907999999999

$ws.Range("H136").Value = 3456.5833
$ws.Range("I136").Value = 3456.5833
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10369.7499
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7819.749899999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7177.0938
$ws.Range("I20").Value = 5384.185
$ws.Range("J20").Value = 16858.8
$ws.Range("K20").Value = 5384.185
$ws.Range("L20").Value = 16858.8
$ws.Range("M20").Value = -5137.185
$ws.Range("N20").Value = -17352.8

$ws.Range("H22").Value = 405.875
$ws.Range("I22").Value = 356.7143
$ws.Range("K22").Value = 356.7143
$ws.Range("M22").Value = -183.7143

$ws.Range("H94").Value = 1219.9
$ws.Range("I94").Value = 1237.9445
$ws.Range("J94").Value = 1057.5
$ws.Range("K94").Value = 1237.9445
$ws.Range("L94").Value = 1057.5
$ws.Range("M94").Value = -786.9445000000001
$ws.Range("N94").Value = -1959.5

$ws.Range("H105").Value = 5557623.5
$ws.Range("I105").Value = 8335477.5
$ws.Range("J105").Value = 1915.8334
$ws.Range("K105").Value = 8335477.5
$ws.Range("L105").Value = 1915.8334
$ws.Range("M105").Value = -8333730.5
$ws.Range("N105").Value = -5409.8334

$ws.Range("H107").Value = 5893.4
$ws.Range("I107").Value = 6616.75
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 6616.75
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -4696.75
$ws.Range("N107").Value = -6840

$ws.Range("H134").Value = 1867.0769
$ws.Range("I134").Value = 1707.091
$ws.Range("J134").Value = 2747
$ws.Range("K134").Value = 5121.272999999999
$ws.Range("L134").Value = 8241
$ws.Range("M134").Value = -2586.272999999999
$ws.Range("N134").Value = -13311

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 19749.5
$ws.Range("I23").Value = 19749.5
$ws.Range("K23").Value = 19749.5
$ws.Range("M23").Value = -19509.5

$ws.Range("H27").Value = 19749.5
$ws.Range("I27").Value = 19749.5
$ws.Range("K27").Value = 19749.5
$ws.Range("M27").Value = -19557.5

$ws.Range("H31").Value = 2050.5
$ws.Range("I31").Value = 1832.909
$ws.Range("K31").Value = 1832.909
$ws.Range("M31").Value = -1537.909

$ws.Range("H34").Value = 2050.5
$ws.Range("I34").Value = 1832.909
$ws.Range("K34").Value = 1832.909
$ws.Range("M34").Value = -1630.909

$ws.Range("H51").Value = 14999.333
$ws.Range("J51").Value = 14999.333
$ws.Range("L51").Value = 14999.333
$ws.Range("N51").Value = -16471.333

$ws.Range("H60").Value = 10237.538
$ws.Range("J60").Value = 10454.091
$ws.Range("L60").Value = 10454.091
$ws.Range("N60").Value = -11476.091

$ws.Range("H61").Value = 14999.333
$ws.Range("J61").Value = 14999.333
$ws.Range("L61").Value = 14999.333
$ws.Range("N61").Value = -15695.333

$ws.Range("H62").Value = 3488.3845
$ws.Range("I62").Value = 2498.111
$ws.Range("J62").Value = 5716.5
$ws.Range("K62").Value = 2498.111
$ws.Range("L62").Value = 5716.5
$ws.Range("M62").Value = -1874.111
$ws.Range("N62").Value = -6964.5

$ws.Range("H65").Value = 3488.3845
$ws.Range("I65").Value = 2498.111
$ws.Range("J65").Value = 5716.5
$ws.Range("K65").Value = 12490.555
$ws.Range("L65").Value = 28582.5
$ws.Range("M65").Value = -9370.555
$ws.Range("N65").Value = -34822.5

$ws.Range("H99").Value = 1999
$ws.Range("I99").Value = 1999
$ws.Range("K99").Value = 1999
$ws.Range("M99").Value = -501

$ws.Range("H107").Value = 1521.9
$ws.Range("J107").Value = 1779
$ws.Range("L107").Value = 1779
$ws.Range("N107").Value = -5619

$ws.Range("H122").Value = 2494.35
$ws.Range("I122").Value = 1816.8823
$ws.Range("K122").Value = 5450.6469
$ws.Range("M122").Value = -3000.6469

$ws.Range("H126").Value = 1999
$ws.Range("I126").Value = 1999
$ws.Range("K126").Value = 5997
$ws.Range("M126").Value = -3527

$ws.Range("H132").Value = 1145
$ws.Range("I132").Value = 1090
$ws.Range("K132").Value = 3270
$ws.Range("M132").Value = -740

$ws.Range("H134").Value = 2698.7144
$ws.Range("I134").Value = 1975.4
$ws.Range("J134").Value = 4507
$ws.Range("K134").Value = 5926.200000000001
$ws.Range("L134").Value = 13521
$ws.Range("M134").Value = -3391.200000000001
$ws.Range("N134").Value = -18591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9020529
$ws.Range("I4").Value = 482118.94
$ws.Range("K4").Value = 1446356.82
$ws.Range("M4").Value = -1446244.82

$ws.Range("H26").Value = 232.125
$ws.Range("I26").Value = 29.333334
$ws.Range("J26").Value = 353.8
$ws.Range("K26").Value = 88.00000199999999
$ws.Range("L26").Value = 1061.4
$ws.Range("M26").Value = 199.999998
$ws.Range("N26").Value = -1637.4

$ws.Range("H75").Value = 3255.6365
$ws.Range("I75").Value = 306.5
$ws.Range("J75").Value = 3911
$ws.Range("K75").Value = 919.5
$ws.Range("L75").Value = 11733
$ws.Range("M75").Value = 78.5
$ws.Range("N75").Value = -13729

$ws.Range("H78").Value = 3255.6365
$ws.Range("I78").Value = 306.5
$ws.Range("J78").Value = 3911
$ws.Range("K78").Value = 2758.5
$ws.Range("L78").Value = 35199
$ws.Range("M78").Value = 2233.5
$ws.Range("N78").Value = -45183

$ws.Range("H82").Value = 4932.6
$ws.Range("J82").Value = 4996
$ws.Range("L82").Value = 14988
$ws.Range("N82").Value = -15800

$ws.Range("H85").Value = 4932.6
$ws.Range("J85").Value = 4996
$ws.Range("L85").Value = 14988
$ws.Range("N85").Value = -17796

$ws.Range("H86").Value = 381.83334
$ws.Range("I86").Value = 378
$ws.Range("K86").Value = 1134
$ws.Range("M86").Value = 52

$ws.Range("H89").Value = 381.83334
$ws.Range("I89").Value = 378
$ws.Range("K89").Value = 3402
$ws.Range("M89").Value = 2526

$ws.Range("H97").Value = 1846.4445
$ws.Range("I97").Value = 925
$ws.Range("J97").Value = 2109.7144
$ws.Range("K97").Value = 2775
$ws.Range("L97").Value = 6329.1432
$ws.Range("M97").Value = -2279
$ws.Range("N97").Value = -7321.1432

$ws.Range("H107").Value = 1203.4286
$ws.Range("I107").Value = 850
$ws.Range("J107").Value = 1286.5883
$ws.Range("K107").Value = 2550
$ws.Range("L107").Value = 3859.7649
$ws.Range("M107").Value = -630
$ws.Range("N107").Value = -7699.7649

$ws.Range("H114").Value = 28572780
$ws.Range("I114").Value = 66667736
$ws.Range("K114").Value = 200003208
$ws.Range("M114").Value = -199999954

$ws.Range("H117").Value = 2688.3333
$ws.Range("J117").Value = 2962
$ws.Range("L117").Value = 8886
$ws.Range("N117").Value = -15770

$ws.Range("H120").Value = 21666.666
$ws.Range("I120").Value = 5000
$ws.Range("J120").Value = 30000
$ws.Range("K120").Value = 15000
$ws.Range("L120").Value = 90000
$ws.Range("M120").Value = -10162
$ws.Range("N120").Value = -99676

$ws.Range("H122").Value = 587.3043
$ws.Range("I122").Value = 333.25
$ws.Range("J122").Value = 640.7895
$ws.Range("K122").Value = 2999.25
$ws.Range("L122").Value = 5767.1055
$ws.Range("M122").Value = -549.25
$ws.Range("N122").Value = -10667.1055

$ws.Range("H129").Value = 3987.8
$ws.Range("I129").Value = 2690
$ws.Range("J129").Value = 4312.25
$ws.Range("K129").Value = 8070
$ws.Range("L129").Value = 12936.75
$ws.Range("M129").Value = -3070
$ws.Range("N129").Value = -22936.75

$ws.Range("H131").Value = 97106.62
$ws.Range("I131").Value = 1595.375
$ws.Range("J131").Value = 155882.77
$ws.Range("K131").Value = 4786.125
$ws.Range("L131").Value = 467648.3099999999
$ws.Range("M131").Value = 253.875
$ws.Range("N131").Value = -477728.3099999999

$ws.Range("H132").Value = 1848.1666
$ws.Range("J132").Value = 2077.8
$ws.Range("L132").Value = 18700.2
$ws.Range("N132").Value = -23760.2

$ws.Range("H136").Value = 6049.0835
$ws.Range("I136").Value = 2933.1667
$ws.Range("J136").Value = 9165
$ws.Range("K136").Value = 8799.500100000001
$ws.Range("L136").Value = 27495
$ws.Range("M136").Value = -3699.500100000001
$ws.Range("N136").Value = -37695

$ws.Range("H139").Value = 2917
$ws.Range("I139").Value = 3365.6667
$ws.Range("J139").Value = 2244
$ws.Range("K139").Value = 10097.0001
$ws.Range("L139").Value = 6732
$ws.Range("M139").Value = -4957.000100000001
$ws.Range("N139").Value = -17012

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 47406
$ws.Range("I26").Value = 23038
$ws.Range("J26").Value = 53498
$ws.Range("K26").Value = 23038
$ws.Range("L26").Value = 53498
$ws.Range("M26").Value = -22758
$ws.Range("N26").Value = -54058

$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 20000
$ws.Range("K29").Value = 20000
$ws.Range("M29").Value = -19710

$ws.Range("H50").Value = 47406
$ws.Range("I50").Value = 23038
$ws.Range("J50").Value = 53498
$ws.Range("K50").Value = 23038
$ws.Range("L50").Value = 53498
$ws.Range("M50").Value = -22540
$ws.Range("N50").Value = -54494

$ws.Range("H97").Value = 30437.545
$ws.Range("I97").Value = 50963.152
$ws.Range("J97").Value = 789.44446
$ws.Range("K97").Value = 50963.152
$ws.Range("L97").Value = 789.44446
$ws.Range("M97").Value = -50467.152
$ws.Range("N97").Value = -1781.44446

$ws.Range("H104").Value = 30671
$ws.Range("J104").Value = 30671
$ws.Range("L104").Value = 30671
$ws.Range("N104").Value = -37659

$ws.Range("H107").Value = 21603.959
$ws.Range("I107").Value = 27335.445
$ws.Range("J107").Value = 4409.5
$ws.Range("K107").Value = 27335.445
$ws.Range("L107").Value = 4409.5
$ws.Range("M107").Value = -25415.445
$ws.Range("N107").Value = -8249.5

$ws.Range("H113").Value = 38465350
$ws.Range("I113").Value = 50003160
$ws.Range("J113").Value = 5993.6665
$ws.Range("K113").Value = 50003160
$ws.Range("L113").Value = 5993.6665
$ws.Range("M113").Value = -50000990
$ws.Range("N113").Value = -10333.6665

$ws.Range("H122").Value = 4327.647
$ws.Range("I122").Value = 4162.1816
$ws.Range("J122").Value = 4631
$ws.Range("K122").Value = 12486.5448
$ws.Range("L122").Value = 13893
$ws.Range("M122").Value = -10036.5448
$ws.Range("N122").Value = -18793

$ws.Range("H126").Value = 3633.3333
$ws.Range("J126").Value = 4800
$ws.Range("L126").Value = 14400
$ws.Range("N126").Value = -19340

$ws.Range("H132").Value = 2049.6667
$ws.Range("I132").Value = 1930
$ws.Range("J132").Value = 2468.5
$ws.Range("K132").Value = 5790
$ws.Range("L132").Value = 7405.5
$ws.Range("M132").Value = -3260
$ws.Range("N132").Value = -12465.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25049.584
$ws.Range("J7").Value = 4060
$ws.Range("L7").Value = 4060
$ws.Range("N7").Value = -4284

$ws.Range("H40").Value = 2939
$ws.Range("I40").Value = 2283.0715
$ws.Range("K40").Value = 2283.0715
$ws.Range("M40").Value = -2147.0715

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H46").Value = 72853.336
$ws.Range("I46").Value = 86424
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 86424
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -86236
$ws.Range("N46").Value = -5376

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

$ws.Range("H74").Value = 25108.5
$ws.Range("I74").Value = 25000
$ws.Range("K74").Value = 25000
$ws.Range("M74").Value = -24002

$ws.Range("H77").Value = 25108.5
$ws.Range("I77").Value = 25000
$ws.Range("K77").Value = 75000
$ws.Range("M77").Value = -70008

$ws.Range("H82").Value = 1325.8889
$ws.Range("I82").Value = 733.75
$ws.Range("J82").Value = 1799.6
$ws.Range("K82").Value = 733.75
$ws.Range("L82").Value = 1799.6
$ws.Range("M82").Value = -372.75
$ws.Range("N82").Value = -2521.6

$ws.Range("H85").Value = 1325.8889
$ws.Range("I85").Value = 733.75
$ws.Range("J85").Value = 1799.6
$ws.Range("K85").Value = 733.75
$ws.Range("L85").Value = 1799.6
$ws.Range("M85").Value = 514.25
$ws.Range("N85").Value = -4295.6

$ws.Range("H100").Value = 136052.4
$ws.Range("I100").Value = 4532.3335
$ws.Range("K100").Value = 4532.3335
$ws.Range("M100").Value = -3991.3335

$ws.Range("H122").Value = 3935.818
$ws.Range("I122").Value = 3555.5715
$ws.Range("J122").Value = 4601.25
$ws.Range("K122").Value = 10666.7145
$ws.Range("L122").Value = 13803.75
$ws.Range("M122").Value = -8216.7145
$ws.Range("N122").Value = -18703.75

$ws.Range("H126").Value = 25049.584
$ws.Range("J126").Value = 4060
$ws.Range("L126").Value = 12180
$ws.Range("N126").Value = -17120

$ws.Range("H132").Value = 10709.154
$ws.Range("I132").Value = 11922.5
$ws.Range("J132").Value = 6664.6665
$ws.Range("K132").Value = 35767.5
$ws.Range("L132").Value = 19993.9995
$ws.Range("M132").Value = -33237.5
$ws.Range("N132").Value = -25053.9995

$ws.Range("H136").Value = 4276.8887
$ws.Range("I136").Value = 2582
$ws.Range("J136").Value = 7666.6665
$ws.Range("K136").Value = 7746
$ws.Range("L136").Value = 22999.9995
$ws.Range("M136").Value = -5196
$ws.Range("N136").Value = -28099.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 2500
$ws.Range("J29").Value = 2500
$ws.Range("L29").Value = 2500
$ws.Range("N29").Value = -3080

$ws.Range("H54").Value = 18767.25
$ws.Range("I54").Value = 8356.666999999999
$ws.Range("J54").Value = 49999
$ws.Range("K54").Value = 8356.666999999999
$ws.Range("L54").Value = 49999
$ws.Range("M54").Value = -7836.666999999999
$ws.Range("N54").Value = -51039

$ws.Range("H81").Value = 2106.4546
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 4000
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 2106.4546
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 20000
$ws.Range("N84").Value = -30608

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 50000680
$ws.Range("I107").Value = 820.6
$ws.Range("J107").Value = 100000536
$ws.Range("K107").Value = 2461.8
$ws.Range("L107").Value = 300001608
$ws.Range("M107").Value = -541.8000000000002
$ws.Range("N107").Value = -300005448

$ws.Range("H113").Value = 1008.8182
$ws.Range("J113").Value = 1333.3334
$ws.Range("L113").Value = 4000.0002
$ws.Range("N113").Value = -8340.0002

$ws.Range("H122").Value = 938.75
$ws.Range("I122").Value = 938.75
$ws.Range("K122").Value = 2816.25
$ws.Range("M122").Value = -366.25

$ws.Range("H136").Value = 1273.1154
$ws.Range("I136").Value = 1113.1305
$ws.Range("K136").Value = 3339.3915
$ws.Range("M136").Value = -789.3914999999997
